# Auto-generated edit script for cryptos.xlsx update
# Applies the cell-value changes described by the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.102.42"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "3.444.82"
$ws.Range("E3").Value = "  +1.99%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "409.28"
$ws.Range("E5").Value = "  +0.97%  "
$ws.Range("D6").Value = "129.25"
$ws.Range("E6").Value = "  -4.67%  "
$ws.Range("D7").Value = "0.622"
$ws.Range("E7").Value = "  +4.89%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("D9").Value = "0.758"
$ws.Range("E9").Value = "  +12.65%  "
$ws.Range("E10").Value = "  +17.59%  "
$ws.Range("D11").Value = "43.18"
$ws.Range("E11").Value = "  +1.79%  "
$ws.Range("E12").Value = "  -0.38%  "
$ws.Range("D13").Value = "8.74"
$ws.Range("E13").Value = "  +4.33%  "
$ws.Range("D14").Value = "20.45"
$ws.Range("E14").Value = "  +3.48%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000199"
$ws.Range("E15").Value = "  +56.72%  "
$ws.Range("D16").Value = "3.430.52"
$ws.Range("E16").Value = "  +1.94%  "
$ws.Range("B17").Value = "Polygon"
$ws.Range("C17").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D17").Value = "1.05"
$ws.Range("E17").Value = "  +2.30%  "
$ws.Range("B18").Value = "WrappedBTC"
$ws.Range("C18").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D18").Value = "62.150.04"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").Value = "11.42"
$ws.Range("E19").Value = "  +1.92%  "
$ws.Range("D20").Value = "367.98"
$ws.Range("E20").Value = "  +20.36%  "
$ws.Range("D21").Value = "87.01"
$ws.Range("E21").Value = "  +4.14%  "
$ws.Range("D22").Value = "3.19"
$ws.Range("E22").Value = "  -1.04%  "
$ws.Range("D23").Value = "13.29"
$ws.Range("E23").Value = "  +4.03%  "
$ws.Range("D24").Value = "3.19"
$ws.Range("E24").Value = "  +1.62%  "
$ws.Range("D25").Value = "31.49"
$ws.Range("E25").Value = "  +6.84%  "
$ws.Range("D26").Value = "4.79"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "8.34"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("E28").Value = "  +2.13%  "
$ws.Range("D29").Value = "2.71"
$ws.Range("E29").Value = "  +9.16%  "
$ws.Range("D30").Value = "44.13"
$ws.Range("E30").Value = "  +5.17%  "
$ws.Range("E31").Value = "  -1.26%  "
$ws.Range("E32").Value = "  -0.53%  "
$ws.Range("D33").Value = "11.79"
$ws.Range("E33").Value = "  +3.74%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").Value = "0.0494"
$ws.Range("E35").Value = "  +3.04%  "
$ws.Range("D36").Value = "51.92"
$ws.Range("E36").Value = "  +0.13%  "
$ws.Range("D37").Value = "0.998"
$ws.Range("E37").Value = "  +0.09%  "
$ws.Range("D38").Value = "3.37"
$ws.Range("E38").Value = "  -1.51%  "
$ws.Range("D39").Value = "2.92"
$ws.Range("E39").Value = "  +0.17%  "
$ws.Range("B40").Value = "TheGraph"
$ws.Range("C40").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D40").Value = "0.315"
$ws.Range("E40").Value = "  +8.39%  "
$ws.Range("B41").Value = "Monero"
$ws.Range("C41").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D41").Value = "142.99"
$ws.Range("E41").Value = "  +3.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.130"
$ws.Range("E42").Value = "  +4.87%  "
$ws.Range("D43").Value = "1.98"
$ws.Range("E43").Value = "  -0.67%  "
$ws.Range("D44").Value = "3.99"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.70"
$ws.Range("E45").Value = "  -0.85%  "
$ws.Range("D46").Value = "2.34"
$ws.Range("E46").Value = "  +4.52%  "
$ws.Range("D47").Value = "21.66"
$ws.Range("E47").Value = "  +1.07%  "
$ws.Range("D48").Value = "2.119.29"
$ws.Range("E48").Value = "  -0.06%  "
$ws.Range("D49").Value = "2.29"
$ws.Range("E49").Value = "  -0.18%  "
$ws.Range("D50").Value = "1.94"
$ws.Range("E50").Value = "  +3.05%  "
$ws.Range("D51").Value = "0.0367"
$ws.Range("E51").Value = "  +8.18%  "
